$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.423.30"
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("D3").Value = "3.781.00"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'626.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.55%  "
$ws.Range("D6").Value = "'166.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").Value = "3.780.43"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("D11").Value = "'0.459"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("E12").Value = "  +2.73%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Value = "4.414.45"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "3.783.32"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "69.436.91"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "'17.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "'470.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("D22").Value = "'9.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").Value = "'0.708"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("E24").Value = "  +5.31%  "
$ws.Range("D25").Value = "'83.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "'12.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.60%  "
$ws.Range("E27").Value = "  +4.79%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "3.930.35"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +4.02%  "
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("D33").Value = "'7.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'28.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'9.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.732.61"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "'0.162"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.26%  "
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("D40").Value = "'3.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.73%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").Value = "'0.969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.299"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("D46").Value = "'43.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Value = "'153.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "'46.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E49").Value = "  +4.95%  "
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("D51").Value = "'1.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
